$d = $word.ActiveDocument

# Locate the "Unidad 1: Algoritmos" text (title line under the main "Reflexión" heading)
# so the edit is anchored to content, not a hard-coded offset.
$find = $d.Content.Duplicate
$find.Find.ClearFormatting()
$find.Find.Text = "Unidad 1: Algoritmos"
$find.Find.Execute() | Out-Null

if (-not $find.Found) {
    throw "Could not find 'Unidad 1: Algoritmos' text to update."
}

# Clear the matched text, collapsing the range to the insertion point while
# leaving the owning paragraph (and its pPr / jc / spacing) untouched.
$find.Text = ""

# Rebuild the line as: "Unidad " (bold) + "2" (bold) + ":" (bold) + " " (plain)
# + "Estructura de control condicional" (plain) -- five separate runs, matching
# how Word splits runs when the unit number/title are edited piecemeal.
# The enclosing <w:p> keeps its original identity/formatting attributes and
# <w:pPr> (centered, no spacing) untouched, same as in the source paragraph.
$xml = '<w:p w14:paraId="59F5F707" w14:textId="2B0CEFAE" w:rsidR="008E05FE" w:rsidRDefault="008E05FE" w:rsidP="001968CB" ' +
       'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' +
       'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
       '<w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/></w:pPr>' +
       '<w:r w:rsidRPr="008E05FE"><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Unidad </w:t></w:r>' +
       '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>2</w:t></w:r>' +
       '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>:</w:t></w:r>' +
       '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
       '<w:r><w:t>Estructura de control condicional</w:t></w:r>' +
       '</w:p>'

$find.InsertXML($xml)
